$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.178.25"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").Value = "2.787.67"
$ws.Range("E3").Value = "  +1.52%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "116.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.53%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.85%  "

$ws.Range("E11").Value = "  +3.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.28%  "

$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("D15").Value = "3.226.50"
$ws.Range("E15").Value = "  +1.67%  "

$ws.Range("D16").Value = "2.807.73"
$ws.Range("E16").Value = "  +2.89%  "

$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Value = "52.065.15"
$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("E19").Value = "  +6.21%  "

$ws.Range("E20").Value = "  +3.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.61%  "

$ws.Range("E22").Value = "  +1.83%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.67%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.11%  "

$ws.Range("E25").Value = "  +6.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.16%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.140"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.33%  "

$ws.Range("E34").Value = "  -0.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0409"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.73%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.12%  "

$ws.Range("E39").Value = "  -1.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +20.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.31%  "

$ws.Range("E44").Value = "  +2.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.20%  "

$ws.Range("D47").Value = "2.075.65"
$ws.Range("E47").Value = "  -1.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.911"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.34%  "

$ws.Range("E51").Value = "  -1.54%  "
